$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: duplicate of row 4's content, with column B replaced ---
$ws.Range("A5").Value = "FF5AD68E11B81949E05377690F803E90"
$ws.Range("B5").Value = "google.com"
$ws.Range("C5").Value = 9357631
$ws.Range("D5").Value = "Infoblox"
$ws.Range("F5").Value = "CNAME"
$ws.Range("G5").Value = "d3bxschxt4niqn.cloudfront.net"
$ws.Range("H5").Value = "AS0018318"
$ws.Range("I5").Value = "Demised"
$ws.Range("K5").Value = 9357631
$ws.Range("L5").Value = "AMH Open Banking"
$ws.Range("M5").Value = "Active"
$ws.Range("N5").Value = $false
$ws.Range("O5").Value = $false
$ws.Range("P5").Value = 9358034
$ws.Range("Q5").Value = "OPEN-BANKING-AMH-HK"
$ws.Range("R5").Value = "IT Service"
$ws.Range("S5").Value = "No"
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = $false
$ws.Range("AF5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AH5").Value = $false
$ws.Range("AI5").Value = $false
$ws.Range("AJ5").Value = $false
$ws.Range("AK5").Value = $false
$ws.Range("AL5").Value = $false
$ws.Range("AM5").Value = $false
$ws.Range("AN5").Value = $false
$ws.Range("AO5").Value = $false
$ws.Range("AP5").Value = "########"
$ws.Range("AQ5").Value = 9358034
$ws.Range("AR5").Value = "OPEN-BANKING-AMH-HK"
$ws.Range("AS5").Value = "Tier 2"
$ws.Range("AT5").Value = "Wholesale Technology"
$ws.Range("AU5").Value = "WS Tech Open Banking"
$ws.Range("AV5").Value = "Retail Banking & Wealth Mgmt"
$ws.Range("AW5").Value = 35032001
$ws.Range("AX5").Value = "Sudipt D SARKAR"
$ws.Range("AY5").Value = "sudipt.d.sarkar@hsbc.co.in"
$ws.Range("AZ5").Value = 4921569
$ws.Range("BA5").Value = "Dominic J PARSONS"
$ws.Range("BB5").Value = "dominicparsons@hsbc.com"
$ws.Range("BC5").Value = 35016742
$ws.Range("BD5").Value = "Jimmy K C MAK"
$ws.Range("BE5").Value = "jimmy.k.c.mak@hsbc.com.hk"
$ws.Range("BF5").Value = $true
$ws.Range("BG5").Value = 50619

# --- Hyperlinks for the new row (mirrors row 4's links) + the new B5 link ---
$ws.Hyperlinks.Add($ws.Range("G5"), "http://d3bxschxt4niqn.cloudfront.net/", "", "", "http://d3bxschxt4niqn.cloudfront.net/")
$ws.Hyperlinks.Add($ws.Range("AY5"), "mailto:sudipt.d.sarkar@hsbc.co.in", "", "", "mailto:sudipt.d.sarkar@hsbc.co.in")
$ws.Hyperlinks.Add($ws.Range("BB5"), "mailto:dominicparsons@hsbc.com", "", "", "mailto:dominicparsons@hsbc.com")
$ws.Hyperlinks.Add($ws.Range("BE5"), "mailto:jimmy.k.c.mak@hsbc.com.hk", "", "", "mailto:jimmy.k.c.mak@hsbc.com.hk")
$ws.Hyperlinks.Add($ws.Range("B5"), "http://www.google.com/", "", "", "www.google.com")

# Re-apply the Hyperlink cell style (Hyperlinks.Add's own auto-style uses a
# different cellXf than the workbook's existing "Hyperlink" named style, so
# restore it explicitly to match the other linked cells in the sheet).
$ws.Range("G5").Style = "Hyperlink"
$ws.Range("AY5").Style = "Hyperlink"
$ws.Range("BB5").Style = "Hyperlink"
$ws.Range("BE5").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 22

# --- Sheet view selection ---
$ws.Range("B19").Select()
